$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(2,2).Value = "C01号直流"
$ws.Cells.Item(2,3).Value = "2025-01-25 13:46:36"
$ws.Cells.Item(2,4).Value = 45949.321608796294

$ws.Cells.Item(3,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(3,2).Value = "C02号直流"
$ws.Cells.Item(3,3).Value = "2025-01-25 17:13:47"
$ws.Cells.Item(3,4).Value = 45949.321608796294

$ws.Cells.Item(4,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(4,2).Value = "C03号直流"
$ws.Cells.Item(4,3).Value = "2025-01-25 14:14:24"
$ws.Cells.Item(4,4).Value = 45949.321608796294

$ws.Cells.Item(5,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(5,2).Value = "C04号直流"
$ws.Cells.Item(5,3).Value = "2025-01-25 06:24:40"
$ws.Cells.Item(5,4).Value = 45949.321608796294

$ws.Cells.Item(6,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(6,2).Value = "C05号直流"
$ws.Cells.Item(6,3).Value = "2025-01-25 16:01:40"
$ws.Cells.Item(6,4).Value = 45949.321608796294

$ws.Cells.Item(7,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(7,2).Value = "D01号直流"
$ws.Cells.Item(7,3).Value = "2025-01-25 18:30:24"
$ws.Cells.Item(7,4).Value = 45949.321608796294

$ws.Cells.Item(8,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(8,2).Value = "D02号直流"
$ws.Cells.Item(8,3).Value = "2025-01-25 15:39:19"
$ws.Cells.Item(8,4).Value = 45949.321608796294

$ws.Cells.Item(9,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(9,2).Value = "D03号直流"
$ws.Cells.Item(9,3).Value = "2025-01-25 16:09:35"
$ws.Cells.Item(9,4).Value = 45949.321608796294

$ws.Cells.Item(10,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10,2).Value = "D04号直流"
$ws.Cells.Item(10,3).Value = "2025-01-25 18:29:02"
$ws.Cells.Item(10,4).Value = 45949.321608796294

$ws.Cells.Item(11,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(11,2).Value = "D05号直流"
$ws.Cells.Item(11,3).Value = "2025-01-25 18:27:29"
$ws.Cells.Item(11,4).Value = 45949.321608796294

$ws.Cells.Item(12,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(12,2).Value = "E01号直流"
$ws.Cells.Item(12,3).Value = "2025-01-25 15:22:58"
$ws.Cells.Item(12,4).Value = 45949.321608796294

$ws.Cells.Item(13,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(13,2).Value = "E02号直流"
$ws.Cells.Item(13,3).Value = "2025-01-25 16:45:57"
$ws.Cells.Item(13,4).Value = 45949.321608796294

$ws.Cells.Item(14,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(14,2).Value = "E03号直流"
$ws.Cells.Item(14,3).Value = "2025-01-25 02:54:59"
$ws.Cells.Item(14,4).Value = 45949.321608796294

$ws.Cells.Item(15,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(15,2).Value = "E04号直流"
$ws.Cells.Item(15,3).Value = "2025-01-25 17:08:37"
$ws.Cells.Item(15,4).Value = 45949.321608796294

$ws.Cells.Item(16,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(16,2).Value = "004B号直流"
$ws.Cells.Item(16,3).Value = "2025-02-19 00:26:27"
$ws.Cells.Item(16,4).Value = 45949.321608796294

$ws.Cells.Item(17,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(17,2).Value = "701号直流"
$ws.Cells.Item(17,3).Value = 45927.457337962966
$ws.Cells.Item(17,4).Value = 45949.321608796294

$ws.Cells.Item(18,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18,2).Value = "804号直流"
$ws.Cells.Item(18,3).Value = 45947.628194444442
$ws.Cells.Item(18,4).Value = 45949.321608796294

$ws.Cells.Item(19,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19,2).Value = "803号直流"
$ws.Cells.Item(19,3).Value = 45948.036365740743
$ws.Cells.Item(19,4).Value = 45949.321608796294

$ws.Cells.Item(20,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(20,2).Value = "501号直流"
$ws.Cells.Item(20,3).Value = 45948.043263888889
$ws.Cells.Item(20,4).Value = 45949.321608796294

$ws.Cells.Item(21,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21,2).Value = "603号直流"
$ws.Cells.Item(21,3).Value = 45948.051631944443
$ws.Cells.Item(21,4).Value = 45949.321608796294

$ws.Cells.Item(22,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22,2).Value = "502号直流"
$ws.Cells.Item(22,3).Value = 45948.208819444444
$ws.Cells.Item(22,4).Value = 45949.321608796294

$ws.Cells.Item(23,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(23,2).Value = "602号直流"
$ws.Cells.Item(23,3).Value = 45948.244259259256
$ws.Cells.Item(23,4).Value = 45949.321608796294

$ws.Cells.Item(24,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(24,2).Value = "B01号直流"
$ws.Cells.Item(24,3).Value = 45948.279756944445
$ws.Cells.Item(24,4).Value = 45949.321608796294

$ws.Cells.Item(25,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(25,2).Value = "801号直流"
$ws.Cells.Item(25,3).Value = 45948.362951388888
$ws.Cells.Item(25,4).Value = 45949.321608796294

$ws.Cells.Item(26,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26,2).Value = "505号直流"
$ws.Cells.Item(26,3).Value = 45948.410497685189
$ws.Cells.Item(26,4).Value = 45949.321608796294

$ws.Cells.Item(27,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(27,2).Value = "208号直流"
$ws.Cells.Item(27,3).Value = 45948.420266203706
$ws.Cells.Item(27,4).Value = 45949.321608796294

$ws.Cells.Item(28,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(28,2).Value = "103号直流"
$ws.Cells.Item(28,3).Value = 45948.458298611113
$ws.Cells.Item(28,4).Value = 45949.321608796294

$ws.Cells.Item(29,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(29,2).Value = "405号直流"
$ws.Cells.Item(29,3).Value = 45948.538113425922
$ws.Cells.Item(29,4).Value = 45949.321608796294

$ws.Cells.Item(30,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30,2).Value = "009B号直流"
$ws.Cells.Item(30,3).Value = 45948.542083333334
$ws.Cells.Item(30,4).Value = 45949.321608796294

$ws.Cells.Item(31,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(31,2).Value = "102号直流"
$ws.Cells.Item(31,3).Value = 45948.545277777775
$ws.Cells.Item(31,4).Value = 45949.321608796294

$ws.Cells.Item(32,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(32,2).Value = "406号直流"
$ws.Cells.Item(32,3).Value = 45948.545682870368
$ws.Cells.Item(32,4).Value = 45949.321608796294

$ws.Cells.Item(33,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(33,2).Value = "111号直流"
$ws.Cells.Item(33,3).Value = 45948.549988425926
$ws.Cells.Item(33,4).Value = 45949.321608796294

$ws.Cells.Item(34,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(34,2).Value = "902号直流"
$ws.Cells.Item(34,3).Value = 45948.563935185186
$ws.Cells.Item(34,4).Value = 45949.321608796294

$ws.Cells.Item(35,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(35,2).Value = "210号直流"
$ws.Cells.Item(35,3).Value = 45948.569212962961
$ws.Cells.Item(35,4).Value = 45949.321608796294

$ws.Cells.Item(36,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(36,2).Value = "405号直流"
$ws.Cells.Item(36,3).Value = 45948.580300925925
$ws.Cells.Item(36,4).Value = 45949.321608796294

$ws.Cells.Item(37,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(37,2).Value = "011A号直流"
$ws.Cells.Item(37,3).Value = 45948.582905092589
$ws.Cells.Item(37,4).Value = 45949.321608796294

$ws.Cells.Item(38,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(38,2).Value = "401号直流"
$ws.Cells.Item(38,3).Value = 45948.597986111112
$ws.Cells.Item(38,4).Value = 45949.321608796294

$ws.Cells.Item(39,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(39,2).Value = "904号直流"
$ws.Cells.Item(39,3).Value = 45948.610115740739
$ws.Cells.Item(39,4).Value = 45949.321608796294

$ws.Cells.Item(40,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(40,2).Value = "504号直流"
$ws.Cells.Item(40,3).Value = 45948.615023148152
$ws.Cells.Item(40,4).Value = 45949.321608796294

$ws.Cells.Item(41,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(41,2).Value = "305号直流"
$ws.Cells.Item(41,3).Value = 45948.616979166669
$ws.Cells.Item(41,4).Value = 45949.321608796294

$ws.Cells.Item(42,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(42,2).Value = "403号直流"
$ws.Cells.Item(42,3).Value = 45948.644699074073
$ws.Cells.Item(42,4).Value = 45949.321608796294

$ws.Cells.Item(43,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(43,2).Value = "006B号直流"
$ws.Cells.Item(43,3).Value = 45948.657407407409
$ws.Cells.Item(43,4).Value = 45949.321608796294

$ws.Cells.Item(44,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(44,2).Value = "106号直流"
$ws.Cells.Item(44,3).Value = 45948.662395833337
$ws.Cells.Item(44,4).Value = 45949.321608796294

$ws.Cells.Item(45,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(45,2).Value = "206号直流"
$ws.Cells.Item(45,3).Value = 45948.662951388891
$ws.Cells.Item(45,4).Value = 45949.321608796294

$ws.Cells.Item(46,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(46,2).Value = "402号直流"
$ws.Cells.Item(46,3).Value = 45948.698460648149
$ws.Cells.Item(46,4).Value = 45949.321608796294

$ws.Cells.Item(47,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(47,2).Value = "A01号直流"
$ws.Cells.Item(47,3).Value = 45948.701874999999
$ws.Cells.Item(47,4).Value = 45949.321608796294

$ws.Cells.Item(48,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(48,2).Value = "203号直流"
$ws.Cells.Item(48,3).Value = 45948.704525462963
$ws.Cells.Item(48,4).Value = 45949.321608796294

$ws.Cells.Item(49,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(49,2).Value = "805号直流"
$ws.Cells.Item(49,3).Value = 45948.77648148148
$ws.Cells.Item(49,4).Value = 45949.321608796294

$ws.Cells.Item(50,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(50,2).Value = "705号直流"
$ws.Cells.Item(50,3).Value = 45948.792719907404
$ws.Cells.Item(50,4).Value = 45949.321608796294

$ws.Cells.Item(51,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(51,2).Value = "604号直流"
$ws.Cells.Item(51,3).Value = 45948.796550925923
$ws.Cells.Item(51,4).Value = 45949.321608796294

$ws.Cells.Item(52,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(52,2).Value = "303号直流"
$ws.Cells.Item(52,3).Value = 45948.802164351851
$ws.Cells.Item(52,4).Value = 45949.321608796294

[void]$ws.Range("G20").Select()

Write-Host "edit complete"